# Update the contact list data on Sheet1 and remove the now-extra last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 21 (name, address, contact)
$data = @(
    @("CERIMONIAL SHOW", "Mais de 10 anos no mercado", "Av. Treze, n° 7A"),
    @("Comitivas - Espaço de Eventos", "R. Comitivas, das Dr. Roberto Cardoso Alves, 5038", "(17) 98145-4054"),
    @("Espaço Magui", "Rua CJ10, R. Vitalina Silva Neves, 38", "(17) 98810-6080"),
    @("Espaço Berrantão", "Rodovia Brigadeiro Faria Lima, km 428, s/n", "Não Informado"),
    @("Camarote Brahma", "Av. dos bretes, s/n", "(17) 3321-0000"),
    @("Barretão Tendas", "Rua Ibrahim Chehadi", "(17) 98108-2227"),
    @("Badalasom - Locação de Geradores Barretos 24h", "(17) 99777-2001", "Não Informado"),
    @("Ramos Vídeo", "Av. Quarenta e Cinco (Numeração Com Zero Inicial), No. 0184", "(17) 98156-2308"),
    @("Nobreza Pulpitos", "Av. Onze, 342", "(17) 98100-0515"),
    @("Parque do Peão Os Independentes", "Rod. Brg. Faria Lima, Km 428", "(17) 3321-0000"),
    @("ESPAÇO DE LAZER BAMBÃ BARRETOS", "Av. Sf- Um", "(17) 99226-7578"),
    @("Doce boutique confeitaria", "R. João Luís da Silva", "(17) 98825-2587"),
    @("Power Técnica", "R. Dez, 187", "(17) 99112-2009"),
    @("Fábrica de Salgados", "Av. Hércules Brasolin, 339", "(17) 3322-7775"),
    @("Os Independentes", "Rod. Brg. Faria Lima, 428", "(17) 3321-0001"),
    @("Recinto Paulo de Lima Correa", "Av. Vinte e Três", "Não Informado"),
    @("Trajes Masculinos a Rigor", "Av. Vinte e Sete, n 442", "(17) 99161-3312"),
    @("Perfumaria Emy", "R. Vinte, 727", "(17) 3325-7474"),
    @("Conka Filmes", "(17) 99256-2426", "Não Informado"),
    @("Machado Alfaiate - Costura e Bordados", "Mais de 15 anos no mercado", "esquina com a - Av. 9, R. Uruguai, 2681")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the now-unused last row (previously row 22) entirely.
$ws.Rows.Item(22).Delete()
